$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.862.69"

$ws.Range("D3").Value = "1.642.52"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  -0.63%  "

$ws.Range("D5").Value = "'216.60"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("E6").Value = "  +1.68%  "

$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").Value = "'0.0622"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").Value = "'19.75"
$ws.Range("E10").Value = "  +4.11%  "

$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.51%  "

$ws.Range("D12").Value = "1.873.01"
$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.12"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.597.50"
$ws.Range("E14").Value = "  -3.31%  "

$ws.Range("E15").Value = "  +1.23%  "

$ws.Range("D16").Value = "'66.02"
$ws.Range("E16").Value = "  +3.05%  "

$ws.Range("D17").Value = "26.891.96"
$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").Value = "'219.16"
$ws.Range("E19").Value = "  +3.71%  "

$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("E21").Value = "  +1.46%  "

$ws.Range("D22").Value = "'6.60"
$ws.Range("E22").Value = "  +6.64%  "

$ws.Range("D23").Value = "'2.40"
$ws.Range("E23").Value = "  +3.61%  "

$ws.Range("D24").Value = "'9.16"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").Value = "'145.97"
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").Value = "'7.40"
$ws.Range("E27").Value = "  +5.92%  "

$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("D29").Value = "'15.77"
$ws.Range("E29").Value = "  +1.57%  "

$ws.Range("D30").Value = "'0.0504"
$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("E33").Value = "  +1.78%  "

$ws.Range("D34").Value = "'1.55"
$ws.Range("E34").Value = "  +2.38%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").Value = "1.246.10"
$ws.Range("E36").Value = "  -1.37%  "

$ws.Range("E37").Value = "  +1.14%  "

$ws.Range("D38").Value = "'0.538"
$ws.Range("E38").Value = "  +2.81%  "

$ws.Range("D39").Value = "'0.829"
$ws.Range("E39").Value = "  +3.36%  "

$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("E41").Value = "  +0.83%  "

$ws.Range("E42").Value = "  +1.85%  "

$ws.Range("D43").Value = "1.785.00"

$ws.Range("D44").Value = "'2.09"
$ws.Range("E44").Value = "  -2.90%  "

$ws.Range("D45").Value = "'60.86"
$ws.Range("E45").Value = "  +1.96%  "

$ws.Range("D46").Value = "'91.28"
$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("E48").Value = "  +15.90%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.61"
$ws.Range("E50").Value = "  +2.08%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0971"
$ws.Range("E51").Value = "  +1.49%  "
